# "final experiments completed with double setup"
# - Row 13 (Code/Description) was updated from the placeholder "C1" /
#   "Capacitive?" entry to the finalized "Q1" / "Capacitive" entry.
# - The sheet's active selection moved from D22 to C14 (and the view
#   scrolled so row 4 is back at the top) to reflect where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Capacitive?" placeholder row to the finalized "Capacitive" entry.
$ws.Range("B13").Value = "Q1"
$ws.Range("C13").Value = "Capacitive"

# Reflect the new selection / scroll position left by the author.
$ws.Range("C14").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
